# Refreshes the cryptos list (Coin / Link / Price / Volume(1h) columns) to
# the latest scraped snapshot, per the commit message "Updated cryptos list
# ... with GitHub Actions". Only the rows whose figures moved are touched;
# two coin pairs (WrappedBTC<->Chainlink, NEARProtocol<->PEPE) swapped rank
# so their Coin/Link/Price/Volume cells are rewritten together.
#
# Price-column values that are plain decimals (e.g. "2.71", "640.50",
# "1.00") are prefixed with a leading apostrophe — Excel's standard way of
# forcing literal text entry — so they stay Text cells (matching the
# original inlineStr cells) instead of being auto-parsed into Number cells
# and losing significant trailing zeros / exact formatting. Values that are
# already unambiguous as text (percent strings, multi-dot thousands values,
# coin names, URLs) are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.407.35'
$ws.Range("E2").Value = '  -2.17%  '
$ws.Range("D3").Value = '3.607.11'
$ws.Range("E3").Value = '  -3.08%  '
$ws.Range("D4").Value = "'2.71"
$ws.Range("E4").Value = '  +25.60%  '
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = "'223.07"
$ws.Range("E6").Value = '  -5.90%  '
$ws.Range("D7").Value = "'640.50"
$ws.Range("E7").Value = '  -2.51%  '
$ws.Range("D8").Value = "'0.420"
$ws.Range("E8").Value = '  -5.11%  '
$ws.Range("E9").Value = '  +5.79%  '
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").Value = '3.603.90'
$ws.Range("E11").Value = '  -3.17%  '
$ws.Range("D12").Value = "'50.83"
$ws.Range("E12").Value = '  +13.48%  '
$ws.Range("E13").Value = '  +4.75%  '
$ws.Range("E14").Value = '  -6.70%  '
$ws.Range("D15").Value = "'6.48"
$ws.Range("E15").Value = '  -5.03%  '
$ws.Range("D16").Value = '4.278.79'
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '95.283.48'
$ws.Range("E17").Value = '  -2.13%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = "'24.51"
$ws.Range("E18").Value = '  +29.18%  '
$ws.Range("D19").Value = "'9.18"
$ws.Range("E19").Value = '  +3.64%  '
$ws.Range("D20").Value = "'13.73"
$ws.Range("E20").Value = '  +5.02%  '
$ws.Range("D21").Value = '3.613.48'
$ws.Range("E21").Value = '  -2.71%  '
$ws.Range("D22").Value = "'0.289"
$ws.Range("E22").Value = '  +35.29%  '
$ws.Range("D23").Value = "'0.531"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").Value = "'137.03"
$ws.Range("E24").Value = '  +16.27%  '
$ws.Range("D25").Value = "'532.66"
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("D26").Value = "'3.27"
$ws.Range("E26").Value = '  -5.48%  '
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").Value = "'7.05"
$ws.Range("E27").Value = '  +2.20%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = "'0.0000202"
$ws.Range("E28").Value = '  -10.19%  '
$ws.Range("D29").Value = "'13.16"
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("D30").Value = '3.774.62'
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("D31").Value = "'13.41"
$ws.Range("E31").Value = '  +5.49%  '
$ws.Range("D32").Value = "'3.13"
$ws.Range("E32").Value = '  +3.25%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = '  +2.56%  '
$ws.Range("D35").Value = "'33.77"
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").Value = "'0.634"
$ws.Range("E36").Value = '  +6.18%  '
$ws.Range("D37").Value = "'0.182"
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").Value = "'0.0554"
$ws.Range("E39").Value = '  +22.03%  '
$ws.Range("D40").Value = "'7.34"
$ws.Range("E40").Value = '  +7.38%  '
$ws.Range("D42").Value = "'8.58"
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("D43").Value = "'590.46"
$ws.Range("E43").Value = '  -7.81%  '
$ws.Range("D44").Value = "'1.02"
$ws.Range("E44").Value = '  +5.87%  '
$ws.Range("D45").Value = "'0.501"
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").Value = "'40.89"
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("D47").Value = "'0.156"
$ws.Range("E47").Value = '  -6.82%  '
$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("D49").Value = "'9.31"
$ws.Range("E49").Value = '  +6.26%  '
$ws.Range("D50").Value = "'231.91"
$ws.Range("E50").Value = '  +11.39%  '
$ws.Range("E51").Value = '  -2.63%  '
